$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 34: RMO No. 27-2021 (Final version for TargetCY 2021)
$ws.Range("B34").Value = "RMO No. 27-2021"
$ws.Range("C34").Value = "October 25, 2021"
$ws.Range("A34").Value = "https://www.bir.gov.ph/images/bir_files/internal_communications_3/Attachments%20of%20RMO%202021/RMO%20No.%2027-2021_Annexes.pdf"
$ws.Range("D34").Value = "2021"
$ws.Range("E34").Value = "Final"
$ws.Range("F34").Value = 2081161
$ws.Range("G34").Formula = "=1015431+39712"
$ws.Range("H34").Value = 305218
$ws.Range("I34").Value = 378721
$ws.Range("J34").Value = 127860
$ws.Range("K34").Value = 214219

# Row 35: RMO No. 16-2022 (Initial version for TargetCY 2022)
$ws.Range("B35").Value = "RMO No. 16-2022"
$ws.Range("C35").Value = "March 15, 2022"
$ws.Range("D35").Value = "2022"
$ws.Range("A35").Value = "https://www.bir.gov.ph/images/bir_files/internal_communications_3/2022/Attachments/RMO16-2022_Annexes.pdf"
$ws.Range("E35").Value = "Initial"
$ws.Range("F35").Value = 2438302
$ws.Range("G35").Value = 1225218.0830000001
$ws.Range("H35").Value = 348345
$ws.Range("I35").Value = 501631.55300000001
$ws.Range("J35").Value = 136741.364
$ws.Range("K35").Value = 226366

# Row 36: RMO No. 30-2022 (Revision 1 for TargetCY 2022)
$ws.Range("C36").Value = "June 9, 2022"
$ws.Range("B36").Value = "RMO No. 30-2022"
$ws.Range("A36").Value = "https://www.bir.gov.ph/images/bir_files/internal_communications_3/2022/Attachments/RMO%2030/RMO%20No.%2030-2022%20Table%205A-F.pdf"
$ws.Range("D36").Value = "2022"
$ws.Range("E36").Value = "Revision 1"
$ws.Range("F36").Value = 2438302
$ws.Range("G36").Value = 1225218.0830000001
$ws.Range("H36").Value = 348345
$ws.Range("I36").Value = 501631.55300000001
$ws.Range("J36").Value = 136741.364
$ws.Range("K36").Value = 226366

# Update the selection to match the post-edit state (next empty row selected)
$null = $ws.Range("A37").Select()
